$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E retain their text formatting so numeric-looking
# strings (e.g. "1.00", "552.80") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "59.890.09"
$ws.Range("E2").Value = "  -3.40%  "

$ws.Range("D3").Value = "3.284.14"
$ws.Range("E3").Value = "  -4.02%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "552.80"
$ws.Range("E5").Value = "  -4.38%  "

$ws.Range("D6").Value = "140.20"
$ws.Range("E6").Value = "  -8.15%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.288.09"
$ws.Range("E8").Value = "  -3.90%  "

$ws.Range("D10").Value = "7.76"
$ws.Range("E10").Value = "  -3.73%  "

$ws.Range("E11").Value = "  -4.94%  "

$ws.Range("D12").Value = "0.405"
$ws.Range("E12").Value = "  -2.87%  "

$ws.Range("D13").Value = "3.849.23"
$ws.Range("E13").Value = "  -3.91%  "

$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("D15").Value = "26.76"
$ws.Range("E15").Value = "  -6.90%  "

$ws.Range("D16").Value = "3.279.56"
$ws.Range("E16").Value = "  -4.24%  "

$ws.Range("E17").Value = "  -5.13%  "

$ws.Range("D18").Value = "59.972.84"
$ws.Range("E18").Value = "  -3.42%  "

$ws.Range("D19").Value = "6.06"
$ws.Range("E19").Value = "  -6.30%  "

$ws.Range("D20").Value = "13.76"
$ws.Range("E20").Value = "  -5.26%  "

$ws.Range("D21").Value = "8.51"
$ws.Range("E21").Value = "  -4.95%  "

$ws.Range("D22").Value = "371.67"
$ws.Range("E22").Value = "  -2.98%  "

$ws.Range("D23").Value = "73.48"
$ws.Range("E23").Value = "  -2.37%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "0.531"
$ws.Range("E25").Value = "  -6.82%  "

$ws.Range("D26").Value = "3.419.93"
$ws.Range("E26").Value = "  -4.01%  "

$ws.Range("E27").Value = "  -9.93%  "

$ws.Range("E28").Value = "  -5.56%  "

$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").Value = "7.06"
$ws.Range("E30").Value = "  -8.45%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  -4.99%  "

$ws.Range("D33").Value = "7.45"
$ws.Range("E33").Value = "  -5.78%  "

$ws.Range("D34").Value = "22.44"
$ws.Range("E34").Value = "  -3.45%  "

$ws.Range("E35").Value = "  -7.39%  "

$ws.Range("D36").Value = "5.05"
$ws.Range("E36").Value = "  -8.76%  "

$ws.Range("D37").Value = "166.17"
$ws.Range("E37").Value = "  -1.38%  "

$ws.Range("D38").Value = "1.51"
$ws.Range("E38").Value = "  -6.86%  "

$ws.Range("D39").Value = "6.62"
$ws.Range("E39").Value = "  -4.46%  "

$ws.Range("D40").Value = "3.317.97"
$ws.Range("E40").Value = "  -3.93%  "

$ws.Range("D41").Value = "26.21"
$ws.Range("E41").Value = "  -16.03%  "

$ws.Range("D42").Value = "0.0724"
$ws.Range("E42").Value = "  -7.55%  "

$ws.Range("D43").Value = "41.64"

$ws.Range("E44").Value = "  -4.32%  "

$ws.Range("D45").Value = "4.10"
$ws.Range("E45").Value = "  -7.04%  "

$ws.Range("D46").Value = "1.56"
$ws.Range("E46").Value = "  -6.91%  "

$ws.Range("E47").Value = "  -6.13%  "

$ws.Range("D49").Value = "2.328.47"
$ws.Range("E49").Value = "  -8.73%  "

$ws.Range("D50").Value = "6.35"
$ws.Range("E50").Value = "  -7.49%  "

$ws.Range("D51").Value = "21.11"
$ws.Range("E51").Value = "  -6.54%  "
